$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-50 down to 12-51.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly observation.
$ws.Range("A11").Value2 = 11
$ws.Range("B11").Value2 = "Vega Monumental Concepción"
$ws.Range("C11").Value2 = "Bíobío"
$ws.Range("D11").Value2 = 45238
$ws.Range("E11").Value2 = 8
$ws.Range("F11").Value2 = 100112022
$ws.Range("G11").Value2 = "Arveja Verde"
$ws.Range("H11").Value2 = "Sin especificar"
$ws.Range("I11").Value2 = "Primera"
$ws.Range("J11").Value2 = 100
$ws.Range("K11").Value2 = 21000
$ws.Range("L11").Value2 = 22000
$ws.Range("M11").Value2 = 21500
$ws.Range("N11").Value2 = "`$/saco 25 kilos"
$ws.Range("O11").Value2 = "Región del Maule"
$ws.Range("P11").Value2 = 860
$ws.Range("Q11").Value2 = 25
$ws.Range("R11").Value2 = "Hortaliza"
